$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Sheet1" to "Results"
$ws.Name = "Results"

# Move the active selection to D349
[void]$ws.Range("D349").Select()
